$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.362.52"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.870.13"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'243.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4708"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("D8").Value = "'0.2877"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "'0.06451"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").Value = "'22.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.871.74"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "'96.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "'0.7245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "'5.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'279.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "30.355.56"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'0.000007497"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "2.110.81"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'5.236"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'6.230"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'163.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "'9.045"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "'0.09639"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "'1.320"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'4.212"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'0.04805"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'0.6882"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").Value = "'2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "'0.01880"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "'2.809"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "'74.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'1.935"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").Value = "'0.4220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'0.8240"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "'100.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").Value = "'9.588"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("D48").Value = "'35.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "'6.954"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "'899.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "'0.05721"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.79%  "
